# Update countries & provincias Spain
# - Update the "Datos actualizados" timestamp text
# - Update the daily COVID numbers for several countries
# - Islas Malvinas / Montserrat swap position in the country list (row 213 <-> 214)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 14:29"

# --- Update per-country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5478009
$ws.Range("C4").Value = 1743
$ws.Range("D4").Value = 2876080
$ws.Range("E4").Value = 2430361
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = 171568

# India (row 6)
$ws.Range("B6").Value = 2530943
$ws.Range("C6").Value = 5721
$ws.Range("D6").Value = 1810079
$ws.Range("E6").Value = 671693
$ws.Range("G6").Value = 37
$ws.Range("H6").Value = 49171

# Alemania (row 22)
$ws.Range("B22").Value = 223780
$ws.Range("C22").Value = 6
$ws.Range("E22").Value = 11941

# Australia (row 72)
$ws.Range("D72").Value = 13634
$ws.Range("E72").Value = 9022

# Sudan (row 84)
$ws.Range("B84").Value = 12211
$ws.Range("C84").Value = 49
$ws.Range("D84").Value = 6340
$ws.Range("E84").Value = 5075
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 796

# Croacia (row 103)
$ws.Range("B103").Value = 6420
$ws.Range("C103").Value = 162
$ws.Range("D103").Value = 5193
$ws.Range("E103").Value = 1062
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 165

# Sri Lanka (row 122)
$ws.Range("B122").Value = 2888
$ws.Range("C122").Value = 2
$ws.Range("E122").Value = 211

# Islandia (row 134)
$ws.Range("B134").Value = 1999
$ws.Range("C134").Value = 16
$ws.Range("D134").Value = 1870
$ws.Range("E134").Value = 119

# Gambia (row 139)
$ws.Range("B139").Value = 1689
$ws.Range("C139").Value = 66
$ws.Range("D139").Value = 347
$ws.Range("E139").Value = 1288
$ws.Range("G139").Value = 4
$ws.Range("H139").Value = 54

# Burkina Faso (row 150)
$ws.Range("B150").Value = 1240
$ws.Range("C150").Value = 2
$ws.Range("D150").Value = 1006
$ws.Range("E150").Value = 180

# Vietnam (row 159)
$ws.Range("D159").Value = 447
$ws.Range("E159").Value = 480

# --- Islas Malvinas / Montserrat: the two entries swap places in the country
#     list, so swap the full data rows 213 and 214 (name + every stat column) ---

$a213 = $ws.Cells.Item(213, 1).Value()
$b213 = $ws.Cells.Item(213, 2).Value()
$c213 = $ws.Cells.Item(213, 3).Value()
$d213 = $ws.Cells.Item(213, 4).Value()
$e213 = $ws.Cells.Item(213, 5).Value()
$f213 = $ws.Cells.Item(213, 6).Value()
$g213 = $ws.Cells.Item(213, 7).Value()
$h213 = $ws.Cells.Item(213, 8).Value()

$a214 = $ws.Cells.Item(214, 1).Value()
$b214 = $ws.Cells.Item(214, 2).Value()
$c214 = $ws.Cells.Item(214, 3).Value()
$d214 = $ws.Cells.Item(214, 4).Value()
$e214 = $ws.Cells.Item(214, 5).Value()
$f214 = $ws.Cells.Item(214, 6).Value()
$g214 = $ws.Cells.Item(214, 7).Value()
$h214 = $ws.Cells.Item(214, 8).Value()

$ws.Cells.Item(213, 1).Value = $a214
$ws.Cells.Item(213, 2).Value = $b214
$ws.Cells.Item(213, 3).Value = $c214
$ws.Cells.Item(213, 4).Value = $d214
$ws.Cells.Item(213, 5).Value = $e214
$ws.Cells.Item(213, 6).Value = $f214
$ws.Cells.Item(213, 7).Value = $g214
$ws.Cells.Item(213, 8).Value = $h214

$ws.Cells.Item(214, 1).Value = $a213
$ws.Cells.Item(214, 2).Value = $b213
$ws.Cells.Item(214, 3).Value = $c213
$ws.Cells.Item(214, 4).Value = $d213
$ws.Cells.Item(214, 5).Value = $e213
$ws.Cells.Item(214, 6).Value = $f213
$ws.Cells.Item(214, 7).Value = $g213
$ws.Cells.Item(214, 8).Value = $h213
